# Regenerate the handback status timestamps, as produced by a fresh CI run
# ("Generate Report for Handback"). The handoff/handback datetimes for the
# 62cbcfb9-... entry (row 3) are bumped forward by ~1 minute; row 5 shares
# the same timestamp values as row 3 in the source data, so it updates too.

$wb = $excel.ActiveWorkbook

# zh-cn sheet
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-22 00:19:58"
$wsZhCn.Range("E5").Value = "2016-03-22 00:19:58"
$wsZhCn.Range("H3").Value = "2016-03-22 00:20:22"
$wsZhCn.Range("H5").Value = "2016-03-22 00:20:22"

# de-de sheet
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-22 00:20:02"
$wsDeDe.Range("E5").Value = "2016-03-22 00:20:02"
$wsDeDe.Range("H3").Value = "2016-03-22 00:20:28"
$wsDeDe.Range("H5").Value = "2016-03-22 00:20:28"
